$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of the data rows (A2:B9) - the raw account numbers and
# style-name values are removed, leaving just the header row and the
# formatted-but-empty rows beneath it (matching rows 10+ which were already empty).
# Column A keeps its number formatting/style (s="5"), so only clear its value;
# column B's per-row style is dropped entirely (matching rows 10+ which have no B cell).
$ws.Range("A2:A9").ClearContents()
$ws.Range("B2:B9").Clear()

# Update the selected range shown in the sheet view.
$ws.Range("A2:XFD9").Select()
